$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole date column (including the header) is reformatted as Text
# instead of a date format FIRST, so that the new values below are
# stored as literal text instead of being re-interpreted as dates.
$ws.Range("C1:C21").NumberFormat = "@"

# Column C values were re-typed as free-form text (no longer real dates).
$ws.Range("C2").Value  = "10 1 16"
$ws.Range("C3").Value  = "9.17.2016"
$ws.Range("C4").Value  = "6/6/17"
$ws.Range("C5").Value  = "6/7/17"
$ws.Range("C6").Value  = "6/8/17"
$ws.Range("C7").Value  = "6/9/17"
$ws.Range("C8").Value  = "6/10/17"
$ws.Range("C9").Value  = "6/11/17"
$ws.Range("C10").Value = "6/12/17"
$ws.Range("C11").Value = "9.17.2016"
$ws.Range("C12").Value = "7.26.2015"
$ws.Range("C13").Value = "7.10.2016"
$ws.Range("C14").Value = "12.29.2015"
$ws.Range("C15").Value = "2.20.2015"
$ws.Range("C16").Value = "6/9/17"
$ws.Range("C17").Value = "6/10/17"
$ws.Range("C18").Value = "6/11/17"
$ws.Range("C19").Value = "6/12/17"
$ws.Range("C20").Value = "6/11/17"
$ws.Range("C21").Value = "6/12/17"

# Column C keeps (approximately) the sheet's default width, recorded
# explicitly now that the column carries its own style.
$ws.Columns("C").ColumnWidth = 10

# Selection left on D6 after the edit.
$ws.Range("D6").Select()
